$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had a weekly Cilantro price history ending at row 212
# (dimension A1:R212). Two newer weekly readings are inserted ahead of the
# final (oldest-dated) row, and every row from 186 down gets the values of
# the reading that (in the new, longer history) now belongs one slot later
# -- i.e. the whole block shifts by one row, with two brand-new rows slotted
# in just above the last (unchanged) historical row. Net effect: insert two
# blank rows at 212:213 (pushing the old row 212 down to row 214, values
# untouched) and then rewrite rows 186-213 with their new contents.

# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg,
# F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificacion

$ws.Rows("212:213").Insert()

$rowsData = @(
    @(186, 3,'Femacal de La Calera','Coquimbo',44474,5,100112040,'Cilantro','Sin especificar','Primera',240,2500,2800,2662,'$/docena de atados (3 kilos)','Provincia de Quillota',887,3,'Hortaliza'),
    @(187, 3,'Femacal de La Calera','Coquimbo',44474,5,100112040,'Cilantro','Sin especificar','Segunda',120,2000,2000,2000,'$/docena de atados (3 kilos)','Provincia de Quillota',667,3,'Hortaliza'),
    @(188, 3,'Femacal de La Calera','Coquimbo',44438,5,100112040,'Cilantro','Sin especificar','Primera',300,3300,3500,3420,'$/docena de atados (3 kilos)','Provincia de Quillota',1140,3,'Hortaliza'),
    @(189, 3,'Femacal de La Calera','Coquimbo',44372,5,100112040,'Cilantro','Sin especificar','Primera',120,3000,3000,3000,'$/docena de atados (3 kilos)','Provincia de Quillota',1000,3,'Hortaliza'),
    @(190, 3,'Femacal de La Calera','Coquimbo',44286,5,100112040,'Cilantro','Sin especificar','Primera',140,3000,3500,3214,'$/docena de atados (3 kilos)','Provincia de Quillota',1071,3,'Hortaliza'),
    @(191, 3,'Femacal de La Calera','Coquimbo',44209,5,100112040,'Cilantro','Sin especificar','Primera',130,2500,2500,2500,'$/docena de atados (3 kilos)','Provincia de Quillota',833,3,'Hortaliza'),
    @(192, 3,'Femacal de La Calera','Coquimbo',44356,5,100112040,'Cilantro','Sin especificar','Primera',120,3500,3500,3500,'$/docena de atados (3 kilos)','Provincia de Quillota',1167,3,'Hortaliza'),
    @(193, 3,'Femacal de La Calera','Coquimbo',44160,5,100112040,'Cilantro','Sin especificar','Primera',120,3500,3500,3500,'$/docena de atados (3 kilos)','Provincia de Quillota',1167,3,'Hortaliza'),
    @(194, 3,'Femacal de La Calera','Coquimbo',44351,5,100112040,'Cilantro','Sin especificar','Primera',260,3000,3300,3173,'$/docena de atados (3 kilos)','Provincia de Quillota',1058,3,'Hortaliza'),
    @(195, 3,'Femacal de La Calera','Coquimbo',44365,5,100112040,'Cilantro','Sin especificar','Primera',340,3300,3500,3406,'$/docena de atados (3 kilos)','Provincia de Quillota',1135,3,'Hortaliza'),
    @(196, 3,'Femacal de La Calera','Coquimbo',44306,5,100112040,'Cilantro','Sin especificar','Primera',130,3500,3500,3500,'$/docena de atados (3 kilos)','Provincia de Quillota',1167,3,'Hortaliza'),
    @(197, 3,'Femacal de La Calera','Coquimbo',44215,5,100112040,'Cilantro','Sin especificar','Primera',130,2500,2500,2500,'$/docena de atados (3 kilos)','Provincia de Quillota',833,3,'Hortaliza'),
    @(198, 3,'Femacal de La Calera','Coquimbo',44175,5,100112040,'Cilantro','Sin especificar','Primera',160,3500,3500,3500,'$/docena de atados (3 kilos)','Provincia de Quillota',1167,3,'Hortaliza'),
    @(199, 3,'Femacal de La Calera','Coquimbo',44461,5,100112040,'Cilantro','Sin especificar','Primera',230,2300,2500,2404,'$/docena de atados (3 kilos)','Provincia de Quillota',801,3,'Hortaliza'),
    @(200, 3,'Femacal de La Calera','Coquimbo',44357,5,100112040,'Cilantro','Sin especificar','Primera',160,3500,3500,3500,'$/docena de atados (3 kilos)','Provincia de Quillota',1167,3,'Hortaliza'),
    @(201, 3,'Femacal de La Calera','Coquimbo',44203,5,100112040,'Cilantro','Sin especificar','Primera',180,2500,2500,2500,'$/docena de atados (3 kilos)','Provincia de Quillota',833,3,'Hortaliza'),
    @(202, 3,'Femacal de La Calera','Coquimbo',44162,5,100112040,'Cilantro','Sin especificar','Primera',50,3500,3500,3500,'$/docena de atados (3 kilos)','Provincia de Quillota',1167,3,'Hortaliza'),
    @(203, 3,'Femacal de La Calera','Coquimbo',44410,5,100112040,'Cilantro','Sin especificar','Primera',140,4000,4300,4129,'$/docena de atados (3 kilos)','Provincia de Quillota',1376,3,'Hortaliza'),
    @(204, 3,'Femacal de La Calera','Coquimbo',44411,5,100112040,'Cilantro','Sin especificar','Primera',120,4000,4000,4000,'$/docena de atados (3 kilos)','Provincia de Quillota',1333,3,'Hortaliza'),
    @(205, 3,'Femacal de La Calera','Coquimbo',44257,5,100112040,'Cilantro','Sin especificar','Primera',160,3000,3000,3000,'$/docena de atados (3 kilos)','Provincia de Quillota',1000,3,'Hortaliza'),
    @(206, 3,'Femacal de La Calera','Coquimbo',44244,5,100112040,'Cilantro','Sin especificar','Primera',80,3000,3000,3000,'$/docena de atados (3 kilos)','Provincia de Quillota',1000,3,'Hortaliza'),
    @(207, 3,'Femacal de La Calera','Coquimbo',44176,5,100112040,'Cilantro','Sin especificar','Primera',140,3000,3300,3129,'$/docena de atados (3 kilos)','Provincia de Quillota',1043,3,'Hortaliza'),
    @(208, 3,'Femacal de La Calera','Coquimbo',44239,5,100112040,'Cilantro','Sin especificar','Primera',120,3000,3000,3000,'$/docena de atados (3 kilos)','Provincia de Quillota',1000,3,'Hortaliza'),
    @(209, 3,'Femacal de La Calera','Coquimbo',44376,5,100112040,'Cilantro','Sin especificar','Primera',270,3000,3300,3167,'$/docena de atados (3 kilos)','Provincia de Quillota',1056,3,'Hortaliza'),
    @(210, 3,'Femacal de La Calera','Coquimbo',44292,5,100112040,'Cilantro','Sin especificar','Primera',125,3500,4000,3760,'$/docena de atados (3 kilos)','Provincia de Quillota',1253,3,'Hortaliza'),
    @(211, 3,'Femacal de La Calera','Coquimbo',44358,5,100112040,'Cilantro','Sin especificar','Primera',160,3500,3500,3500,'$/docena de atados (3 kilos)','Provincia de Quillota',1167,3,'Hortaliza'),
    @(212, 3,'Femacal de La Calera','Coquimbo',44211,5,100112040,'Cilantro','Sin especificar','Primera',80,2500,2500,2500,'$/docena de atados (3 kilos)','Provincia de Quillota',833,3,'Hortaliza'),
    @(213, 3,'Femacal de La Calera','Coquimbo',44425,5,100112040,'Cilantro','Sin especificar','Primera',90,3000,3000,3000,'$/docena de atados (3 kilos)','Provincia de Quillota',1000,3,'Hortaliza')
)

foreach ($row in $rowsData) {
    $r = $row[0]
    for ($i = 1; $i -le 18; $i++) {
        $ws.Cells.Item($r, $i).Value = $row[$i]
    }
}
